# Insert a new "Run 50" column of raw results, relabel the old "Mean"
# column (AZ) as "Run 50" holding the new run's raw value, and append a
# fresh "Mean" column (BA) holding the recomputed mean across all 51 runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New run-50 raw values (what used to live in the old "Mean" column, AZ)
$newRunValue = 43.28022961
# Recomputed mean across all 51 runs (old 50-run mean + new run, averaged)
$newMeanValue = 43.34792798

# AZ1 currently reads "Mean" -> relabel it "Run 50"
$ws.Range("AZ1").Value = "Run 50"
# BA1 becomes the new trailing "Mean" header; copy AZ1's header formatting
# (bold/centered/bordered) onto it, then set its text.
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)
$ws.Range("BA1").Value = "Mean"

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 52).Value = $newRunValue   # column AZ
    $ws.Cells.Item($row, 53).Value = $newMeanValue  # column BA
}
